$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E12').Value = 'رقم الموديل 05...'
$ws.Range('E13').Value = 'رقم الموديل 05...'

$ws.Range('E14').Value = 'رقم الموديل A3...'
$ws.Range('E15').Value = 'رقم الموديل A3...'

$ws.Range('E16').Value = 'بلد الصنع اسبا...'

$ws.Range('E18').Value = 'رقم الموديل A8...'

$ws.Range('E47').Value = 'المقاس 90x3 ال...'

$ws.Range('E48').Value = 'المقاس 63x2 ال...'

$ws.Range('E49').Value = 'المقاس 50x11/2...'

$ws.Range('E50').Value = 'المقاس 40x11/4...'

$ws.Range('E51').Value = 'المقاس 20x3/4 ...'

$ws.Range('E62').Value = 'تفاصيل المنتج ...'
$ws.Range('E63').Value = 'تفاصيل المنتج ...'
$ws.Range('E64').Value = 'تفاصيل المنتج ...'
$ws.Range('E65').Value = 'تفاصيل المنتج ...'
$ws.Range('E66').Value = 'تفاصيل المنتج ...'
$ws.Range('E67').Value = 'تفاصيل المنتج ...'
$ws.Range('E108').Value = 'تفاصيل المنتج ...'

$ws.Range('E68').Value = 'المقاس 120*60 ...'
$ws.Range('E69').Value = 'المقاس 120*60 ...'
$ws.Range('E70').Value = 'المقاس 120*60 ...'

$ws.Range('E71').Value = 'بلد الصنع ايطا...'

$ws.Range('E72').Value = 'المقاس 135*300...'

$ws.Range('E73').Value = 'المقاس 240*120...'

$ws.Range('E74').Value = 'مقاس اللوح 20(...'

$ws.Range('E82').Value = 'المقاس 6بوصه ت...'
$ws.Range('E83').Value = 'المقاس 6بوصه ت...'

$ws.Range('E84').Value = 'المقاس 70*70 ت...'

$ws.Range('E85').Value = 'المقاس 100*100...'

$ws.Range('E86').Value = 'المقاس 120*120...'

$ws.Range('E87').Value = 'المقاس 50*50 ت...'

$ws.Range('E89').Value = 'المقاس 3 متر ت...'
$ws.Range('E90').Value = 'المقاس 3 متر ت...'
$ws.Range('E91').Value = 'المقاس 3 متر ت...'

$ws.Range('E92').Value = 'اللون رصاصي تف...'
$ws.Range('E93').Value = 'اللون رصاصي تف...'
$ws.Range('E94').Value = 'اللون رصاصي تف...'

$ws.Range('E95').Value = 'اللون أسود تفا...'
$ws.Range('E96').Value = 'اللون أسود تفا...'
$ws.Range('E97').Value = 'اللون أسود تفا...'
$ws.Range('E98').Value = 'اللون أسود تفا...'
$ws.Range('E99').Value = 'اللون أسود تفا...'
$ws.Range('E100').Value = 'اللون أسود تفا...'

$ws.Range('E101').Value = 'اللون رصاصي*اب...'

$ws.Range('F44').Value = 'المقاس 20مم 25...'

$ws.Range('F45').Value = 'المقاس 20x1/2 ...'
$ws.Range('F46').Value = 'المقاس 20x1/2 ...'

$ws.Range('F108').Value = 'العرض ملم 600م...'

$ws.Range('H14').Value = 'كرسي من ديبا ك...'

$ws.Range('H15').Value = 'كرسي ذا جاب كر...'

$ws.Range('H16').Value = 'مغسلة ذا جاب ا...'

$ws.Range('H32').Value = 'الالوان: SILVE...'

$ws.Range('H33').Value = 'الالوان: BRONZ...'

$ws.Range('H34').Value = 'الالوان: SN ال...'
$ws.Range('H35').Value = 'الالوان: SN ال...'
$ws.Range('H36').Value = 'الالوان: SN ال...'
$ws.Range('H37').Value = 'الالوان: SN ال...'
$ws.Range('H38').Value = 'الالوان: SN ال...'
$ws.Range('H39').Value = 'الالوان: SN ال...'
$ws.Range('H40').Value = 'الالوان: SN ال...'
$ws.Range('H41').Value = 'الالوان: SN ال...'

$ws.Range('H52').Value = 'Gateway G2 هي ...'

$ws.Range('H53').Value = 'تعريف المنتج ج...'

$ws.Range('H54').Value = 'الالوان: MSN ,...'
$ws.Range('H56').Value = 'الالوان: MSN ,...'

$ws.Range('H59').Value = 'الالوان: فضي -...'
$ws.Range('H60').Value = 'الالوان: فضي -...'

$ws.Range('H61').Value = 'الالوان: أسود ...'

$ws.Range('H102').Value = 'خشب بلوط أحمر ...'

$ws.Range('H103').Value = 'خشب الكرز يتحو...'

$ws.Range('H105').Value = 'خشب ميرانتي وع...'

$ws.Range('H107').Value = 'خشب الواوا Waw...'

$ws.Range('H108').Value = 'نظره عامة يتم ...'

$ws.Range('H109').Value = 'شبك مزارع سياج...'

$ws.Range('H110').Value = 'Item Name Text...'

$ws.Range('H111').Value = 'مواسير حديد مو...'
